$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '30.412.53'
$ws.Range("E2").Value = '  +1.24%  '
$ws.Range("D3").Value = '2.001.47'
$ws.Range("E3").Value = '  +4.42%  '
$ws.Range("D4").Value = '1.002'
$ws.Range("E4").Value = '  +0.08%  '
$ws.Range("D5").Value = '324.40'
$ws.Range("E5").Value = '  +1.34%  '
$ws.Range("E6").Value = '  +0.00%  '
$ws.Range("D7").Value = '0.5100'
$ws.Range("E7").Value = '  +1.48%  '
$ws.Range("E8").Value = '  +2.46%  '
$ws.Range("D9").Value = '0.08725'
$ws.Range("E9").Value = '  +5.97%  '
$ws.Range("D10").Value = '1.134'
$ws.Range("E10").Value = '  +2.32%  '
$ws.Range("D11").Value = '43.03'
$ws.Range("E11").Value = '  +2.37%  '
$ws.Range("D12").Value = '24.49'
$ws.Range("E12").Value = '  +3.15%  '
$ws.Range("D13").Value = '1.999.53'
$ws.Range("E13").Value = '  +4.37%  '
$ws.Range("D14").Value = '6.562'
$ws.Range("E14").Value = '  +2.32%  '
$ws.Range("D15").Value = '7.457'
$ws.Range("E15").Value = '  +2.33%  '
$ws.Range("D17").Value = '94.22'
$ws.Range("E17").Value = '  +2.14%  '
$ws.Range("D18").Value = '0.00001115'
$ws.Range("E18").Value = '  +1.67%  '
$ws.Range("D19").Value = '0.06505'
$ws.Range("E19").Value = '  +0.09%  '
$ws.Range("D20").Value = '18.89'
$ws.Range("E21").Value = '  +0.07%  '
$ws.Range("D22").Value = '6.190'
$ws.Range("E22").Value = '  +4.25%  '
$ws.Range("D23").Value = '30.464.17'
$ws.Range("E23").Value = '  +1.25%  '
$ws.Range("D24").Value = '11.91'
$ws.Range("E24").Value = '  +5.65%  '
$ws.Range("D25").Value = '2.225'
$ws.Range("E25").Value = '  +1.46%  '
$ws.Range("D26").Value = '2.230.62'
$ws.Range("E26").Value = '  +4.69%  '
$ws.Range("D27").Value = '22.32'
$ws.Range("E27").Value = '  -0.09%  '
$ws.Range("D28").Value = '163.17'
$ws.Range("E28").Value = '  +0.62%  '
$ws.Range("D29").Value = '2.400'
$ws.Range("E29").Value = '  +4.61%  '
$ws.Range("D30").Value = '131.09'
$ws.Range("E30").Value = '  +1.68%  '
$ws.Range("D31").Value = '1.134'
$ws.Range("E31").Value = '  +0.16%  '
$ws.Range("D32").Value = '0.1052'
$ws.Range("E32").Value = '  +1.27%  '
$ws.Range("D33").Value = '6.063'
$ws.Range("E33").Value = '  +0.75%  '
$ws.Range("D34").Value = '3.848'
$ws.Range("E34").Value = '  +0.83%  '
$ws.Range("D35").Value = '1.340'
$ws.Range("E35").Value = '  +11.03%  '
$ws.Range("E36").Value = '  +3.12%  '
$ws.Range("D37").Value = '5.435'
$ws.Range("E37").Value = '  +1.57%  '
$ws.Range("D38").Value = '0.06602'
$ws.Range("E38").Value = '  +2.63%  '
$ws.Range("D39").Value = '12.49'
$ws.Range("E39").Value = '  +9.55%  '
$ws.Range("D40").Value = '0.2195'
$ws.Range("E40").Value = '  +1.45%  '
$ws.Range("D41").Value = '9.024'
$ws.Range("E41").Value = '  +1.31%  '
$ws.Range("D42").Value = '0.6622'
$ws.Range("E42").Value = '  +2.94%  '
$ws.Range("D43").Value = '1.232'
$ws.Range("E43").Value = '  +1.15%  '
$ws.Range("D44").Value = '13.56'
$ws.Range("E44").Value = '  +1.62%  '
$ws.Range("D45").Value = '0.6157'
$ws.Range("E45").Value = '  +2.66%  '
$ws.Range("D46").Value = '2.190'
$ws.Range("E46").Value = '  -0.16%  '
$ws.Range("D47").Value = '3.665'
$ws.Range("E47").Value = '  +0.77%  '
$ws.Range("D48").Value = '1.265'
$ws.Range("E48").Value = '  +4.13%  '
$ws.Range("D49").Value = '124.47'
$ws.Range("E49").Value = '  +0.80%  '
$ws.Range("D50").Value = '80.38'
$ws.Range("E50").Value = '  +1.92%  '
$ws.Range("D51").Value = '0.06888'
$ws.Range("E51").Value = '  +1.36%  '
